$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1571
$ws.Range("F5").Value = 8915
$ws.Range("F6").Value = 242
$ws.Range("F7").Value = 105
$ws.Range("F8").Value = 1250
$ws.Range("F10").Value = 546
$ws.Range("F14").Value = 284
$ws.Range("F17").Value = 1440
$ws.Range("F18").Value = 1311
$ws.Range("F21").Value = 1344
$ws.Range("F22").Value = 72
$ws.Range("F23").Value = 217
$ws.Range("F26").Value = 38
$ws.Range("F27").Value = 61
$ws.Range("F28").Value = 289
$ws.Range("F29").Value = 1059
$ws.Range("F32").Value = 216
$ws.Range("F33").Value = 181
$ws.Range("F34").Value = 55
$ws.Range("F35").Value = 575
$ws.Range("F36").Value = 602
$ws.Range("F38").Value = 122
$ws.Range("F39").Value = 74
$ws.Range("F40").Value = 148
$ws.Range("F41").Value = 475
$ws.Range("F43").Value = 677
$ws.Range("F44").Value = 201

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 33
$ws.Range("F6").Value = 47
$ws.Range("F16").Value = 660
$ws.Range("F23").Value = 924
$ws.Range("F25").Value = 1028
$ws.Range("F26").Value = 211
$ws.Range("F29").Value = 110
$ws.Range("F31").Value = 144
$ws.Range("F36").Value = 26
$ws.Range("F39").Value = 93

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 738
$ws.Range("F6").Value = 277
$ws.Range("F7").Value = 137
$ws.Range("F8").Value = 2006
$ws.Range("F9").Value = 3014

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1571
$ws.Range("F4").Value = 738
$ws.Range("F6").Value = 8915
$ws.Range("F7").Value = 277
$ws.Range("F8").Value = 137
$ws.Range("F10").Value = 3014
$ws.Range("F12").Value = 105
$ws.Range("F13").Value = 1250
$ws.Range("F17").Value = 284
$ws.Range("F18").Value = 1440
$ws.Range("F19").Value = 1311
$ws.Range("F22").Value = 1344
$ws.Range("F23").Value = 217
$ws.Range("F25").Value = 289
$ws.Range("F26").Value = 289
$ws.Range("F27").Value = 1059
$ws.Range("F30").Value = 924
$ws.Range("F31").Value = 216
$ws.Range("F33").Value = 211
$ws.Range("F34").Value = 55
$ws.Range("F35").Value = 602
$ws.Range("F37").Value = 122
$ws.Range("F38").Value = 74
$ws.Range("F39").Value = 148
$ws.Range("F40").Value = 144
$ws.Range("F41").Value = 475
$ws.Range("F42").Value = 677
$ws.Range("F44").Value = 26
$ws.Range("F45").Value = 201
$ws.Range("F47").Value = 93
